# Money Pig review — headline/meta rewrite + refreshed "What we like" bullets.
#
# Plain Find/Replace on these paragraphs would coalesce the pre-existing
# empty <w:r/> placeholder run into the text run it precedes (since both
# have no distinguishing rPr), which the target revision does not do - it
# only ever edits <w:t> contents. So instead we locate each paragraph by
# its current text and rebuild it via Range.InsertXML, explicitly keeping
# the leading empty run (and any bold/italic run formatting) untouched.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphByText($oldText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.TrimEnd("`r`a") -eq $oldText) {
            return $p
        }
    }
    throw "Paragraph with text '$oldText' not found"
}

function Replace-ParagraphXml($oldText, $innerXml) {
    $p = Find-ParagraphByText $oldText
    $p.Range.InsertXML("<w:p $w>$innerXml</w:p>") | Out-Null
}

# 1. Main H1 title
Replace-ParagraphXml "Play Money Pig Free Slot Game | Review" `
    '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Money Pig for Free - Online Slot Game Review</w:t></w:r>'

# "What we like" bullet list (ListBullet style, leading empty run kept)
$bulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

Replace-ParagraphXml "Customization options for coin value and winning lines" `
    "$bulletPPr<w:r/><w:r><w:t>Customization options for winning lines and coin value</w:t></w:r>"

Replace-ParagraphXml "Gold and purple color scheme" `
    "$bulletPPr<w:r/><w:r><w:t>Captivating atmospheric music</w:t></w:r>"

Replace-ParagraphXml "Simple and modern symbols" `
    "$bulletPPr<w:r/><w:r><w:t>Modern and visually appealing symbols</w:t></w:r>"

Replace-ParagraphXml "Two special symbols (Wild and Scatter)" `
    "$bulletPPr<w:r/><w:r><w:t>Bonus mode with random winnings</w:t></w:r>"

# Bold recap title near the bottom of the page
Replace-ParagraphXml "Play Money Pig Free Slot Game | Review" `
    '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Money Pig for Free - Online Slot Game Review</w:t></w:r>'

# Italic meta description
Replace-ParagraphXml "Read our expert review of Money Pig, an online slot game by Capecod. Play for free and enjoy bonus features, gold and purple colors, and modern symbols." `
    '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Money Pig, a captivating online slot game with customization options. Play for free now!</w:t></w:r>'
